$wb = $excel.ActiveWorkbook

# Rename the first sheet ("practical_ex") to "ExerciseAbstract" - this sheet
# now backs the new ExerciseAbstract DB table (Dao/Repository/ViewModel work
# described in the commit message).
$wsAbstract = $wb.Worksheets.Item("practical_ex")
$wsAbstract.Name = "ExerciseAbstract"

# Move the cursor/selection on "Table_arrangement_2" (previously parked at
# L7) to C9, reflecting where the author was last working in that sheet.
$wsTable2 = $wb.Worksheets.Item("Table_arrangement_2")
$wsTable2.Activate()
$wsTable2.Range("C9").Select()
